$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.133.46"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "1.899.11"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.98"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5231"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3798"
$ws.Range("E8").Value = "  +0.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.30"
$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9054"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08199"
$ws.Range("E12").Value = "  -2.13%  "

$ws.Range("D13").Value = "1.879.89"
$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.48"
$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.349"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("D20").Value = "27.178.53"
$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.119"
$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("D22").Value = "2.120.11"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.329"
$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "149.50"
$ws.Range("E26").Value = "  +2.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.27"
$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.737"
$ws.Range("E28").Value = "  -1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.34"
$ws.Range("E29").Value = "  +0.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.825"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.866"
$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09233"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05048"
$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7928"
$ws.Range("E34").Value = "  -3.25%  "

$ws.Range("E35").Value = "  -0.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.977"
$ws.Range("E36").Value = "  +1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.383"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.652"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5733"
$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("E40").Value = "  +1.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.081"
$ws.Range("E41").Value = "  +0.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.017"
$ws.Range("E42").Value = "  +1.17%  "

$ws.Range("E43").Value = "  -0.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.36"
$ws.Range("E44").Value = "  -1.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1516"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4898"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.16"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "38.57"
$ws.Range("E50").Value = "  +3.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.14"
$ws.Range("E51").Value = "  +0.95%  "
